# "before switch fuel and area": add a GAS_HEATING column (C) to the
# base_data_fuel sheet, mirroring the existing ED_ELEC_APPLIANCES column (B)
# but driven off the new 270176.53 GWh figure instead of 82805.6 GWh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("C1").Value = "GAS_HEATING"

# Formulas mirroring column B, using the new base figure (270176.53)
$ws.Range("C2").Formula = "=270176.53*0.893964"
$ws.Range("C3").Formula = "= 270176.53*8.646003"
$ws.Range("C4").Formula = "= 270176.53*86.46003"

# Match the new column width used for column C (13.42578125 char units)
$ws.Columns.Item(3).ColumnWidth = 12.666666666666666

# Move the active selection to C8, as in the edited workbook
$ws.Range("C8").Select() | Out-Null
